# AFA 2020.xlsx update
# - Fill in results (E/H/K/N/Q/T) for contests 14-18 (rows 23-27)
# - Add 5 new contests (23-27: "RR vs DC", "KXI vs KKR", "CSK vs RCB", "SRH vs RR", "MI vs DC")
#   by inserting 5 rows before the old row 32 "buffer" row, so the totals section
#   (and its SUM ranges) shift down and auto-expand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: fill in match results for contests 14-18 (rows 23-27) ----
$results = @{
    23 = @{ E=100; H=60;  K=80;  N=40; Q=20; T=0  }
    24 = @{ E=60;  H=80;  K=100; N=40; Q=0;  T=20 }
    25 = @{ E=40;  H=100; K=0;   N=60; Q=80; T=20 }
    26 = @{ E=100; H=0;   K=60;  N=80; Q=20; T=40 }
    27 = @{ E=100; H=80;  K=60;  N=20; Q=40; T=0  }
}

foreach ($row in $results.Keys) {
    $vals = $results[$row]
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("T$row").Value = $vals.T
}

# ---- Step 2: insert 5 rows at row 32 so the summary SUM ranges (which end at
#      row 32) expand automatically, and copy formatting from row 31 into them ----
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows("32:32").Insert()
}

$ws.Rows("31:31").Copy()
$ws.Range("A32:U36").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---- Step 3: populate the 5 new contest rows (23-27) ----
$newContests = @{
    32 = @{ Num=23; Match="RR vs DC" }
    33 = @{ Num=24; Match="KXI vs KKR" }
    34 = @{ Num=25; Match="CSK vs RCB" }
    35 = @{ Num=26; Match="SRH vs RR" }
    36 = @{ Num=27; Match="MI vs DC" }
}

foreach ($row in $newContests.Keys) {
    $info = $newContests[$row]
    $ws.Range("A$row").Value = $info.Num
    $ws.Range("B$row").Value = 1
    $ws.Range("C$row").Value = $info.Match

    $ws.Range("D$row").Formula = "=IF(ISERROR(VLOOKUP(RANK(E$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(E$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE))"
    $ws.Range("G$row").Formula = "=IF(ISERROR(VLOOKUP(RANK(H$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(H$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE))"
    $ws.Range("J$row").Formula = "=IF(ISERROR(VLOOKUP(RANK(K$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(K$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE))"
    $ws.Range("M$row").Formula = "=IF(ISERROR(VLOOKUP(RANK(N$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(N$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE))"
    $ws.Range("P$row").Formula = "=IF(ISERROR(VLOOKUP(RANK(Q$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(Q$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE))"
    $ws.Range("S$row").Formula = "=IF(ISERROR(VLOOKUP(RANK(T$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(T$row, (`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  score, 2, FALSE))"
}

# ---- Step 4: select the final totals cell, like the source workbook does ----
$ws.Range("U41").Select()

Write-Host "Done."
